$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1875.875
$ws.Range("I43").Value = 1703.4
$ws.Range("J43").Value = 2163.3333
$ws.Range("K43").Value = 1703.4
$ws.Range("L43").Value = 2163.3333
$ws.Range("M43").Value = -1634.4
$ws.Range("N43").Value = -2301.3333
$ws.Range("H62").Value = 2869.5715
$ws.Range("I62").Value = 2869.5715
$ws.Range("K62").Value = 2869.5715
$ws.Range("M62").Value = -2245.5715
$ws.Range("H65").Value = 2869.5715
$ws.Range("I65").Value = 2869.5715
$ws.Range("K65").Value = 14347.8575
$ws.Range("M65").Value = -11227.8575
$ws.Range("H116").Value = 6604.4546
$ws.Range("I116").Value = 5849.8887
$ws.Range("K116").Value = 5849.8887
$ws.Range("M116").Value = -2407.8887
$ws.Range("H132").Value = 4342.5835
$ws.Range("I132").Value = 4342.5835
$ws.Range("K132").Value = 13027.7505
$ws.Range("M132").Value = -10497.7505
$ws.Range("H137").Value = 3129275.5
$ws.Range("I137").Value = 5002509.5
$ws.Range("K137").Value = 15007528.5
$ws.Range("M137").Value = -15004978.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2460
$ws.Range("I2").Value = 2460
$ws.Range("K2").Value = 2460
$ws.Range("M2").Value = -2347
$ws.Range("H32").Value = 1835.9149
$ws.Range("I32").Value = 1896.4
$ws.Range("K32").Value = 1896.4
$ws.Range("M32").Value = -1609.4
$ws.Range("H45").Value = 16205.689
$ws.Range("I45").Value = 22561.1
$ws.Range("J45").Value = 2082.5557
$ws.Range("K45").Value = 22561.1
$ws.Range("L45").Value = 2082.5557
$ws.Range("M45").Value = -22184.1
$ws.Range("N45").Value = -2836.5557
$ws.Range("H61").Value = 2117.5833
$ws.Range("I61").Value = 1663.2069
$ws.Range("K61").Value = 1663.2069
$ws.Range("M61").Value = -1451.2069
$ws.Range("H102").Value = 3464.8635
$ws.Range("I102").Value = 2619.4119
$ws.Range("J102").Value = 6339.4
$ws.Range("K102").Value = 2619.4119
$ws.Range("L102").Value = 6339.4
$ws.Range("M102").Value = -997.4119000000001
$ws.Range("N102").Value = -9583.4
$ws.Range("H116").Value = 2460
$ws.Range("I116").Value = 2460
$ws.Range("K116").Value = 2460
$ws.Range("M116").Value = -166
$ws.Range("H122").Value = 2915.5
$ws.Range("I122").Value = 2249
$ws.Range("J122").Value = 4248.5
$ws.Range("K122").Value = 6747
$ws.Range("L122").Value = 12745.5
$ws.Range("M122").Value = -4297
$ws.Range("N122").Value = -17645.5
$ws.Range("H132").Value = 2881.6956
$ws.Range("I132").Value = 2634.0588
$ws.Range("J132").Value = 3583.3333
$ws.Range("K132").Value = 7902.176399999999
$ws.Range("L132").Value = 10749.9999
$ws.Range("M132").Value = -5372.176399999999
$ws.Range("N132").Value = -15809.9999
$ws.Range("H136").Value = 2117.5833
$ws.Range("I136").Value = 1663.2069
$ws.Range("K136").Value = 4989.620699999999
$ws.Range("M136").Value = -2439.620699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2460
$ws.Range("I3").Value = 2460
$ws.Range("K3").Value = 2460
$ws.Range("M3").Value = -2346
$ws.Range("H20").Value = 50009300
$ws.Range("I20").Value = 71440000
$ws.Range("K20").Value = 71440000
$ws.Range("M20").Value = -71439753
$ws.Range("H105").Value = 37145450
$ws.Range("I105").Value = 5001400
$ws.Range("J105").Value = 50003068
$ws.Range("K105").Value = 5001400
$ws.Range("L105").Value = 50003068
$ws.Range("M105").Value = -4999653
$ws.Range("N105").Value = -50006562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3381773.2
$ws.Range("I31").Value = 2363.476
$ws.Range("K31").Value = 2363.476
$ws.Range("M31").Value = -2068.476
$ws.Range("H34").Value = 3381773.2
$ws.Range("I34").Value = 2363.476
$ws.Range("K34").Value = 2363.476
$ws.Range("M34").Value = -2161.476
$ws.Range("H105").Value = 2478.2144
$ws.Range("I105").Value = 1836.875
$ws.Range("J105").Value = 3333.3333
$ws.Range("K105").Value = 1836.875
$ws.Range("L105").Value = 3333.3333
$ws.Range("M105").Value = -89.875
$ws.Range("N105").Value = -6827.3333
$ws.Range("H122").Value = 2206.2144
$ws.Range("I122").Value = 2268.2307
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 6804.6921
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -4354.6921
$ws.Range("N122").Value = -9100
$ws.Range("H134").Value = 2362.658
$ws.Range("J134").Value = 2629.375
$ws.Range("L134").Value = 7888.125
$ws.Range("N134").Value = -12958.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 895.9048
$ws.Range("I5").Value = 500.6
$ws.Range("J5").Value = 1884.1666
$ws.Range("K5").Value = 1501.8
$ws.Range("L5").Value = 5652.4998
$ws.Range("M5").Value = -1389.8
$ws.Range("N5").Value = -5876.4998
$ws.Range("H8").Value = 290
$ws.Range("I8").Value = 290
$ws.Range("K8").Value = 870
$ws.Range("M8").Value = -731
$ws.Range("H98").Value = 497.64706
$ws.Range("I98").Value = 334
$ws.Range("J98").Value = 548
$ws.Range("K98").Value = 1002
$ws.Range("L98").Value = 1644
$ws.Range("M98").Value = 496
$ws.Range("N98").Value = -4640
$ws.Range("H122").Value = 1560.1904
$ws.Range("J122").Value = 1697.3684
$ws.Range("L122").Value = 15276.3156
$ws.Range("N122").Value = -20176.3156
$ws.Range("H132").Value = 3340
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -41060
$ws.Range("H135").Value = 895.9048
$ws.Range("I135").Value = 500.6
$ws.Range("J135").Value = 1884.1666
$ws.Range("K135").Value = 4505.400000000001
$ws.Range("L135").Value = 16957.4994
$ws.Range("M135").Value = -1970.400000000001
$ws.Range("N135").Value = -22027.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 5634.385
$ws.Range("I107").Value = 340.72726
$ws.Range("K107").Value = 340.72726
$ws.Range("M107").Value = 1579.27274
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4250.4165
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 14001.667
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 14001.667
$ws.Range("N3").Value = -14225.667
$ws.Range("M3").Value = -888
$ws.Range("H15").Value = 4250.4165
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 14001.667
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 14001.667
$ws.Range("N15").Value = -14341.667
$ws.Range("M15").Value = -830
$ws.Range("H46").Value = 4688.1
$ws.Range("I46").Value = 3554.4285
$ws.Range("K46").Value = 3554.4285
$ws.Range("M46").Value = -3366.4285
$ws.Range("H55").Value = 1106.5834
$ws.Range("I55").Value = 908.7778
$ws.Range("J55").Value = 1700
$ws.Range("K55").Value = 908.7778
$ws.Range("L55").Value = 1700
$ws.Range("M55").Value = -735.7778
$ws.Range("N55").Value = -2046
$ws.Range("H61").Value = 1755.0303
$ws.Range("I61").Value = 1655.6207
$ws.Range("K61").Value = 1655.6207
$ws.Range("M61").Value = -1453.6207
$ws.Range("H113").Value = 1755.0303
$ws.Range("I113").Value = 1655.6207
$ws.Range("K113").Value = 1655.6207
$ws.Range("M113").Value = 514.3793000000001
$ws.Range("H122").Value = 13069.429
$ws.Range("I122").Value = 6250
$ws.Range("J122").Value = 15797.2
$ws.Range("K122").Value = 18750
$ws.Range("L122").Value = 47391.60000000001
$ws.Range("M122").Value = -16300
$ws.Range("N122").Value = -52291.60000000001
$ws.Range("H132").Value = 3665.4062
$ws.Range("I132").Value = 3280.6155
$ws.Range("K132").Value = 9841.8465
$ws.Range("M132").Value = -7311.8465
$ws.Range("H140").Value = 138610
$ws.Range("J140").Value = 138610
$ws.Range("L140").Value = 138610
$ws.Range("N140").Value = -148970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 474747500
$ws.Range("I100").Value = 90909090
$ws.Range("J100").Value = 666666700
$ws.Range("K100").Value = 181818180
$ws.Range("L100").Value = 1333333400
$ws.Range("M100").Value = -181817639
$ws.Range("N100").Value = -1333334482
$ws.Range("H122").Value = 20834904
$ws.Range("I122").Value = 1874.625
$ws.Range("J122").Value = 62500960
$ws.Range("K122").Value = 5623.875
$ws.Range("L122").Value = 187502880
$ws.Range("M122").Value = -3173.875
$ws.Range("N122").Value = -187507780
